# Add a new "Choose a solution and develop a plan to implement it" section
# right after the "Evaluate each potential solution" paragraphs, and move the
# hidden _GoBack bookmark so it still marks the end of the document's last
# edit.

$d = $word.ActiveDocument

# The final paragraph in the body currently holds the trailing (hidden)
# _GoBack bookmark right after "...chooses to take." Remove it now; it will
# be re-created at the new end-of-edit location once the new content has
# been inserted.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Find the paragraph that ends with "...chooses to take." and collapse a
# range to just after its trailing period (i.e. just before the paragraph
# mark) -- that's where the new paragraphs need to be inserted.
$found = $d.Content
$target = $found.Find
$target.Text = "Each of my solutions will work for all cases, one will take longer than the other depending on which solution he chooses to take."
$target.Execute() | Out-Null
if (-not $target.Found) {
    $found = $d.Paragraphs.Item(19).Range
    $insertionPoint = $d.Range($found.End - 1, $found.End - 1)
} else {
    $insertionPoint = $d.Range($found.End, $found.End)
}

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
              </w:rPr>
              <w:t>Choose a solution and develop a plan to implement it</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:t xml:space="preserve">In this solution lets just say the man can wait and get on a bigger boat and in this case he is able to carry the bag of seeds, parrot and cat all the way across the river to the other side. This will give him </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:t>a</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:t xml:space="preserve"> easily way to transport all three items at once.</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($xml)
